# Update betting-odds values in row 4 (Correcaminos vs Atl. Morelia)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value  = 3.05
$ws.Range("I4").Value  = 2.57
$ws.Range("J4").Value  = 3.2
$ws.Range("K4").Value  = 2.07
$ws.Range("L4").Value  = 3.05
$ws.Range("M4").Value  = 8.199999999999999
$ws.Range("O4").Value  = 1.3

$ws.Range("AK4").Value = 30

$ws.Range("AN4").Value = 4.65
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 20
$ws.Range("AQ4").Value = 60
$ws.Range("AR4").Value = 90
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.62
$ws.Range("AU4").Value = 6.4
$ws.Range("AV4").Value = 50
$ws.Range("AW4").Value = 4.55
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 18.5
$ws.Range("AZ4").Value = 55
$ws.Range("BA4").Value = 75
$ws.Range("BB4").Value = 200
